$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("patch riri")

# --- Update I8 / J8 ---
# I8 used to hold "128,138,148,158,168" (the "128 series") and is now replaced by a brand
# new value, the "127 series": "127,137,147,157,167".
# J8 used to hold "129,139,149,159,169" (the "129 series", now retired) and takes on the
# old "128 series" text that used to live in I8: "128,138,148,158,168".
#
# Both strings look like thousand-separated numbers (groups of 3 digits), so a plain
# `.Value = "..."` assignment gets auto-coerced to a number by Excel's smart entry parsing
# (stripping the commas) and also reformats the cell (new style). To store them as literal
# text - exactly like the source file has them (shared string, original style untouched) -
# write the text through a formula and then flatten it to a value via copy / paste-special
# values-only, which does not re-run the numeric auto-detection.

$i8 = $ws.Cells.Item(8, 9)
$i8.Formula = "=""127,137,147,157,167"""
$i8.Copy()
$i8.PasteSpecial(-4163) | Out-Null   # xlPasteValues

$j8 = $ws.Cells.Item(8, 10)
$j8.Formula = "=""128,138,148,158,168"""
$j8.Copy()
$j8.PasteSpecial(-4163) | Out-Null   # xlPasteValues

$excel.CutCopyMode = 0

# --- Update the active selection on the "patch riri" sheet ---
[void]$ws.Range("J8").Select()
